$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 15:52"

# --- Row 32: Arabia Saudita - Casos criticos (F32) ---
$ws.Range("F32").Value = 78

# --- Row 35: Noruega - Casos totales / Nuevos casos / Recuperados ---
$ws.Range("B35").Value = 7036
$ws.Range("C35").Value = 99
$ws.Range("E35").Value = 6841

# --- Row 42: Serbia - Casos criticos (F42) ---
$ws.Range("F42").Value = 126

# --- Row 55: Argentina - Casos activos / Recuperados / Casos criticos ---
$ws.Range("D55").Value = 685
$ws.Range("E55").Value = 1944
$ws.Range("F55").Value = 127

# --- Rows 144/145: Togo and Bermudas swap order + Togo data update ---
# Row 144 becomes Togo with updated stats
$ws.Range("A144").Value = "Togo"
$ws.Range("B144").Value = 84
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 49
$ws.Range("E144").Value = 30
$ws.Range("F144").Value = 0

# Row 145 becomes Bermudas (previous Togo row data moves down, Bermudas keeps its stats)
$ws.Range("A145").Value = "Bermudas"
$ws.Range("B145").Value = 83
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 35
$ws.Range("E145").Value = 43
$ws.Range("F145").Value = 9
